$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accuracy")

$ws.Range("B2").Value = 0.2118068965517241
$ws.Range("B3").Value = 0.1805458229957766
$ws.Range("B4").Value = 0.2521492170709241
$ws.Range("B5").Value = 0.1827738103183608
$ws.Range("B6").Value = 0.2104925837758652
$ws.Range("B7").Value = 0.2757242757242757
$ws.Range("B8").Value = 0.3210372229192807
$ws.Range("B9").Value = 0.271122659259923
$ws.Range("B10").Value = 0.216893039049236
$ws.Range("B11").Value = 0.2738805263656158
$ws.Range("B12").Value = 0.2485025026667761
$ws.Range("B13").Value = 0.2939434534301629
$ws.Range("B14").Value = 0.2832591683289857
$ws.Range("B15").Value = 0.2672127950068266
$ws.Range("B16").Value = 0.2971181376820615
$ws.Range("B17").Value = 0.2026161062297681
